$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Broadmeadows BonBon Bakery (new content, was Craigieburn Line train)
$ws.Range("B4").Value = "BonBon Bakery  Shop G134, Broadmeadows Central  1099/1168 Pascoe Vale Road  Broadmeadows VIC 3047"
$ws.Range("C4").Value = "12:30pm - 12:45pm 9/2/2021"
$ws.Range("D4").Value = "Case attended venue"

# Row 5 - Craigieburn Line train (moved down from row 4, with punctuation fix)
$ws.Range("B5").Value = "Craigieburn Line train"
$ws.Range("C5").Value = "1:25pm - 1:59pm  9/02/2021"
$ws.Range("D5").Value = "Case caught train from Broadmeadows Railway Station to Glenroy Railway Station"

# Row 6 - Woolworths spacing update + time colon fix
$ws.Range("B6").Value = "Woolworths  Broadmeadows Central  Pascoe Vale Road  Broadmeadows VIC 3047"
$ws.Range("C6").Value = "12:15pm - 12:30 pm 9/2/2021"

# Row 9 - Commonwealth Bank comma removed
$ws.Range("B9").Value = "Commonwealth Bank  28-32 Kingsway  Glen Waverley VIC 3015"

# Row 10 - HSBC Bank extra space after comma
$ws.Range("B10").Value = "HSBC Bank,  38 Kingsway  Glen Waverley VIC 3015"

# Row 11 - time colon fix
$ws.Range("C11").Value = "1:35pm  2:17pm  9/02/2021"

# Row 13 - time colon fix
$ws.Range("C13").Value = "6:45am - 7:30am  8/02/21"

# Row 24 - time colon fix
$ws.Range("C24").Value = "4pm - 7:30pm 10/2/2021"
